$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44; existing rows 44-61 shift down to 45-62.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly record.
$ws.Cells.Item(44, 1).Value = 6
$ws.Cells.Item(44, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(44, 3).Value = "Metropolitana"
$ws.Cells.Item(44, 4).Value = (Get-Date -Year 2023 -Month 9 -Day 7 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(44, 5).Value = 13
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100102
$ws.Cells.Item(44, 8).Value = "Cítricos"
$ws.Cells.Item(44, 9).Value = 100102006
$ws.Cells.Item(44, 10).Value = "Pomelo"
$ws.Cells.Item(44, 11).Value = "Start Ruby"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 15
$ws.Cells.Item(44, 14).Value = 190000
$ws.Cells.Item(44, 15).Value = 190000
$ws.Cells.Item(44, 16).Value = 190000
$ws.Cells.Item(44, 17).Value = "`$/bins (350 kilos)"
$ws.Cells.Item(44, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(44, 19).Value = 543
$ws.Cells.Item(44, 20).Value = 350
